$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.743.92"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "2.559.00"
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "514.05"
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("D6").Value = "137.41"
$ws.Range("E6").Value = "  -5.59%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.556"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").Value = "2.578.74"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").Value = "0.0984"
$ws.Range("E11").Value = "  -5.76%  "
$ws.Range("D12").Value = "0.323"
$ws.Range("E12").Value = "  -4.40%  "
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "3.020.53"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "57.706.21"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "19.95"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("D17").Value = "2.591.20"
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("D18").Value = "0.0000130"
$ws.Range("E18").Value = "  -4.93%  "
$ws.Range("D19").Value = "331.08"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").Value = "4.26"
$ws.Range("E20").Value = "  -4.46%  "
$ws.Range("D21").Value = "9.99"
$ws.Range("E21").Value = "  -6.17%  "
$ws.Range("D22").Value = "6.30"
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "65.48"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "0.163"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "0.396"
$ws.Range("E27").Value = "  -5.16%  "
$ws.Range("D28").Value = "6.87"
$ws.Range("E28").Value = "  -4.81%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "0.0₃0702"
$ws.Range("E30").Value = "  -12.58%  "
$ws.Range("E31").Value = "  -8.13%  "
$ws.Range("E32").Value = "  -4.09%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "148.98"
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "18.48"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").Value = "3.86"
$ws.Range("E35").Value = "  -7.40%  "
$ws.Range("D36").Value = "1.10"
$ws.Range("E36").Value = "  -7.98%  "
$ws.Range("D37").Value = "36.06"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").Value = "0.824"
$ws.Range("E38").Value = "  -5.48%  "
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  -5.25%  "
$ws.Range("D40").Value = "1.40"
$ws.Range("E40").Value = "  -5.67%  "
$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  -5.12%  "
$ws.Range("D43").Value = "272.14"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "10.70"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "0.588"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").Value = "0.0935"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.969.04"
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "0.0510"
$ws.Range("E48").Value = "  -4.86%  "
$ws.Range("D49").Value = "18.15"
$ws.Range("E49").Value = "  -6.81%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0216"
$ws.Range("E50").Value = "  -6.11%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "17.66"
$ws.Range("E51").Value = "  -6.00%  "
